# update scripts wuth new tpm
#
# Rewrites the NATMI ligand/receptor summary table (Tnfsf13 -> Tnfrsf11b)
# with values recomputed from the new TPM-based expression matrix.
# The "Target cluster" set expanded from {FAPs, MuSCs} to {ECs, FAPs, MuSCs}
# for every "Sending cluster", so the data block grows from 8 to 12 rows
# (sheet rows 2-13) while columns A:T keep the same layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the data rows (2-13) with the updated TPM-derived values.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric NATMI metrics.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf13"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.342908333333333
$ws.Cells.Item(2, 8).Value = 4.028725
$ws.Cells.Item(2, 9).Value = 0.2879023314891748
$ws.Cells.Item(2, 10).Value = 0.2879023314891748
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.06861733333333334
$ws.Cells.Item(2, 14).Value = 0.205852
$ws.Cells.Item(2, 15).Value = 0.01654048691795588
$ws.Cells.Item(2, 16).Value = 0.01654048691795588
$ws.Cells.Item(2, 17).Value = 0.09214678874444443
$ws.Cells.Item(2, 18).Value = 0.8293210986999999
$ws.Cells.Item(2, 19).Value = 0.004762044747645694
$ws.Cells.Item(2, 20).Value = 0.004762044747645694

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf13"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.342908333333333
$ws.Cells.Item(3, 8).Value = 4.028725
$ws.Cells.Item(3, 9).Value = 0.2879023314891748
$ws.Cells.Item(3, 10).Value = 0.2879023314891748
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.776574666666666
$ws.Cells.Item(3, 14).Value = 11.329724
$ws.Cells.Item(3, 15).Value = 0.9103586635352137
$ws.Cells.Item(3, 16).Value = 0.9103586635352137
$ws.Cells.Item(3, 17).Value = 5.071593591322221
$ws.Cells.Item(3, 18).Value = 45.64434232189999
$ws.Cells.Item(3, 19).Value = 0.2620943817231572
$ws.Cells.Item(3, 20).Value = 0.2620943817231572

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf13"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.342908333333333
$ws.Cells.Item(4, 8).Value = 4.028725
$ws.Cells.Item(4, 9).Value = 0.2879023314891748
$ws.Cells.Item(4, 10).Value = 0.2879023314891748
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.303255
$ws.Cells.Item(4, 14).Value = 0.909765
$ws.Cells.Item(4, 15).Value = 0.07310084954683041
$ws.Cells.Item(4, 16).Value = 0.07310084954683042
$ws.Cells.Item(4, 17).Value = 0.4072436666249999
$ws.Cells.Item(4, 18).Value = 3.665192999625
$ws.Cells.Item(4, 19).Value = 0.02104590501837186
$ws.Cells.Item(4, 20).Value = 0.02104590501837186

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnfsf13"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.4963216666666666
$ws.Cells.Item(5, 8).Value = 1.488965
$ws.Cells.Item(5, 9).Value = 0.1064050028249084
$ws.Cells.Item(5, 10).Value = 0.1064050028249084
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.06861733333333334
$ws.Cells.Item(5, 14).Value = 0.205852
$ws.Cells.Item(5, 15).Value = 0.01654048691795588
$ws.Cells.Item(5, 16).Value = 0.01654048691795588
$ws.Cells.Item(5, 17).Value = 0.03405626924222222
$ws.Cells.Item(5, 18).Value = 0.30650642318
$ws.Cells.Item(5, 19).Value = 0.001759990557230456
$ws.Cells.Item(5, 20).Value = 0.001759990557230456

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnfsf13"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.4963216666666666
$ws.Cells.Item(6, 8).Value = 1.488965
$ws.Cells.Item(6, 9).Value = 0.1064050028249084
$ws.Cells.Item(6, 10).Value = 0.1064050028249084
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.776574666666666
$ws.Cells.Item(6, 14).Value = 11.329724
$ws.Cells.Item(6, 15).Value = 0.9103586635352137
$ws.Cells.Item(6, 16).Value = 0.9103586635352137
$ws.Cells.Item(6, 17).Value = 1.874395832851111
$ws.Cells.Item(6, 18).Value = 16.86956249565999
$ws.Cells.Item(6, 19).Value = 0.09686671616514425
$ws.Cells.Item(6, 20).Value = 0.09686671616514425

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnfsf13"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.4963216666666666
$ws.Cells.Item(7, 8).Value = 1.488965
$ws.Cells.Item(7, 9).Value = 0.1064050028249084
$ws.Cells.Item(7, 10).Value = 0.1064050028249084
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.303255
$ws.Cells.Item(7, 14).Value = 0.909765
$ws.Cells.Item(7, 15).Value = 0.07310084954683041
$ws.Cells.Item(7, 16).Value = 0.07310084954683042
$ws.Cells.Item(7, 17).Value = 0.150512027025
$ws.Cells.Item(7, 18).Value = 1.354608243225
$ws.Cells.Item(7, 19).Value = 0.007778296102533694
$ws.Cells.Item(7, 20).Value = 0.007778296102533696

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Tnfsf13"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.1501973333333333
$ws.Cells.Item(8, 8).Value = 0.450592
$ws.Cells.Item(8, 9).Value = 0.03220038283833477
$ws.Cells.Item(8, 10).Value = 0.03220038283833477
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06861733333333334
$ws.Cells.Item(8, 14).Value = 0.205852
$ws.Cells.Item(8, 15).Value = 0.01654048691795588
$ws.Cells.Item(8, 16).Value = 0.01654048691795588
$ws.Cells.Item(8, 17).Value = 0.01030614048711111
$ws.Cells.Item(8, 18).Value = 0.092755264384
$ws.Cells.Item(8, 19).Value = 0.0005326100110906475
$ws.Cells.Item(8, 20).Value = 0.0005326100110906475

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Tnfsf13"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.1501973333333333
$ws.Cells.Item(9, 8).Value = 0.450592
$ws.Cells.Item(9, 9).Value = 0.03220038283833477
$ws.Cells.Item(9, 10).Value = 0.03220038283833477
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.776574666666666
$ws.Cells.Item(9, 14).Value = 11.329724
$ws.Cells.Item(9, 15).Value = 0.9103586635352137
$ws.Cells.Item(9, 16).Value = 0.9103586635352137
$ws.Cells.Item(9, 17).Value = 0.5672314440675554
$ws.Cells.Item(9, 18).Value = 5.105082996607999
$ws.Cells.Item(9, 19).Value = 0.02931389748602867
$ws.Cells.Item(9, 20).Value = 0.02931389748602867

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Tnfsf13"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1501973333333333
$ws.Cells.Item(10, 8).Value = 0.450592
$ws.Cells.Item(10, 9).Value = 0.03220038283833477
$ws.Cells.Item(10, 10).Value = 0.03220038283833477
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.303255
$ws.Cells.Item(10, 14).Value = 0.909765
$ws.Cells.Item(10, 15).Value = 0.07310084954683041
$ws.Cells.Item(10, 16).Value = 0.07310084954683042
$ws.Cells.Item(10, 17).Value = 0.04554809232
$ws.Cells.Item(10, 18).Value = 0.40993283088
$ws.Cells.Item(10, 19).Value = 0.00235387534121545
$ws.Cells.Item(10, 20).Value = 0.00235387534121545

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Tnfsf13"
$ws.Cells.Item(11, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.675030666666667
$ws.Cells.Item(11, 8).Value = 8.025092000000001
$ws.Cells.Item(11, 9).Value = 0.573492282847582
$ws.Cells.Item(11, 10).Value = 0.573492282847582
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.06861733333333334
$ws.Cells.Item(11, 14).Value = 0.205852
$ws.Cells.Item(11, 15).Value = 0.01654048691795588
$ws.Cells.Item(11, 16).Value = 0.01654048691795588
$ws.Cells.Item(11, 17).Value = 0.1835534709315556
$ws.Cells.Item(11, 18).Value = 1.651981238384
$ws.Cells.Item(11, 19).Value = 0.009485841601989087
$ws.Cells.Item(11, 20).Value = 0.009485841601989087

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Tnfsf13"
$ws.Cells.Item(12, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.675030666666667
$ws.Cells.Item(12, 8).Value = 8.025092000000001
$ws.Cells.Item(12, 9).Value = 0.573492282847582
$ws.Cells.Item(12, 10).Value = 0.573492282847582
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.776574666666666
$ws.Cells.Item(12, 14).Value = 11.329724
$ws.Cells.Item(12, 15).Value = 0.9103586635352137
$ws.Cells.Item(12, 16).Value = 0.9103586635352137
$ws.Cells.Item(12, 17).Value = 10.10245304828978
$ws.Cells.Item(12, 18).Value = 90.922077434608
$ws.Cells.Item(12, 19).Value = 0.5220836681608835
$ws.Cells.Item(12, 20).Value = 0.5220836681608835

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Tnfsf13"
$ws.Cells.Item(13, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.675030666666667
$ws.Cells.Item(13, 8).Value = 8.025092000000001
$ws.Cells.Item(13, 9).Value = 0.573492282847582
$ws.Cells.Item(13, 10).Value = 0.573492282847582
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.303255
$ws.Cells.Item(13, 14).Value = 0.909765
$ws.Cells.Item(13, 15).Value = 0.07310084954683041
$ws.Cells.Item(13, 16).Value = 0.07310084954683042
$ws.Cells.Item(13, 17).Value = 0.81121642482
$ws.Cells.Item(13, 18).Value = 7.300947823380001
$ws.Cells.Item(13, 19).Value = 0.0419227730847094
$ws.Cells.Item(13, 20).Value = 0.04192277308470941
